$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 660 (the "「開⇔閉」" post) entirely, shifting all subsequent rows up by one.
$ws.Rows.Item(660).Delete()
